$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting all existing data down by one row
$ws.Rows.Item(1).Insert()

# Populate the new header row with "source" / "target" labels
$ws.Range("A1").Value = "source"
$ws.Range("B1").Value = "target"

# Update the active selection to A2, matching the recorded view state
$ws.Range("A2").Select()
